$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        if ($val -match ",") {
            $parts = $val -split ","
            for ($i = 0; $i -lt $parts.Length; $i++) {
                $parts[$i] = $parts[$i].Trim()
            }
            $n = $parts.Length
            $reversed = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
